$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that should be stored as a shared *text* string
# (matches how the author's Transaction-date / Local-Payday columns already
# contain text like "20201014" rather than numbers), without leaving a
# lingering explicit cell style behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Helper: write a numeric amount using the same number format already used
# by the existing "Amount" column (style index 1 / #,##0.00).
function Set-AmountValue($range, $value) {
    $range.Value = $value
    $range.NumberFormat = "#,##0.00"
}

# --- Fill in the "Transaction date" helper column (E) for a few existing
#     rows that previously only had the Local Payday (F) populated ---
Set-TextValue $ws.Range("E9") "20201014"
Set-TextValue $ws.Range("E15") "20201014"
Set-TextValue $ws.Range("E19") "20201014"
Set-TextValue $ws.Range("E21") "20201014"

# --- New row 22 ---
Set-TextValue $ws.Range("B22") "20201030"
Set-TextValue $ws.Range("C22") "20201014-ZSAC-0004"
Set-TextValue $ws.Range("D22") "CINTAS CORPORATION"
Set-TextValue $ws.Range("E22") "20201014"
Set-TextValue $ws.Range("F22") "20201014"
Set-AmountValue $ws.Range("G22") 1221.37

# --- New row 23 ---
Set-TextValue $ws.Range("B23") "20201030"
Set-TextValue $ws.Range("C23") "20201014-ZSAC-0005"
Set-TextValue $ws.Range("D23") "CALIFORNIA DEPARTMENT OF TAX AND FEE ADMINISTRATION"
Set-TextValue $ws.Range("F23") "20201023"
Set-AmountValue $ws.Range("G23") 541

# --- New row 24 ---
Set-TextValue $ws.Range("B24") "20201030"
Set-TextValue $ws.Range("C24") "20201014-ZSAC-0006"
Set-TextValue $ws.Range("D24") "MULTI SERVICE AVIATION"
Set-TextValue $ws.Range("F24") "20201021"
Set-AmountValue $ws.Range("G24") 7017.59

# --- Row 26: a single formatted-but-empty helper cell (as in the target
#     workbook), which also pushes the sheet dimension down to row 26 ---
$ws.Range("J26").NumberFormat = "#,##0.00"

# --- Column D got wider (and best-fit) once the longer vendor name was
#     added; set it as close as this runtime's width grid allows ---
$ws.Columns("D").ColumnWidth = 55.6

# --- Update the autofilter to cover the newly added rows. Toggle
#     AutoFilterMode off first so re-applying AutoFilter resizes the
#     existing filter range instead of just removing it ---
$ws.AutoFilterMode = $false
$ws.Range("A1:G26").AutoFilter()

# --- The hidden, sheet-local _FilterDatabase defined name also needs to
#     track the new filter range ---
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$G`$26"

# --- Restore the active selection to match the reviewed workbook ---
$ws.Range("I22").Select()
